# Refresh the crypto price/volume table to the latest scraped values.
# (GitHub Actions scheduled update - mirrors the commit "Updated cryptos
# list ... with GitHub Actions".)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.408.09"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "1.887.76"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.77%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.687"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.358"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0751"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.73%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.775"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.07%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.159.73"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "1.907.67"
$ws.Range("D18").Value = "35.421.37"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").Value = "0.0₃0829"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "245.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.18%  "
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.858"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0735"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("E41").Value = "  +3.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("D45").Value = "1.311.56"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0805"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.77%  "
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.92%  "
